$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'57.141.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Formula = "'2.986.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  -2.18%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'501.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -4.56%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Formula = "'138.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Formula = "'  -3.11%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E8").Formula = "'  -4.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Formula = "'  -4.45%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Formula = "'  -4.43%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Formula = "'0.358"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Formula = "'  -3.79%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'3.499.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Formula = "'  -2.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Formula = "'  -3.64%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Formula = "'  -5.96%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Formula = "'57.172.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Formula = "'  -0.97%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Formula = "'  -3.58%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'2.988.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  -2.35%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Formula = "'  -3.45%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Formula = "'  -3.65%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Formula = "'321.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Formula = "'  -5.07%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Formula = "'  +0.00%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Formula = "'  +0.48%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Formula = "'  -1.90%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Formula = "'63.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Formula = "'  -2.82%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Formula = "'  +0.37%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Formula = "'  -5.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Formula = "'0.0₃0896"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Formula = "'  -8.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Formula = "'6.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Formula = "'7.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  -3.18%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Formula = "'  -4.26%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Formula = "'1.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Formula = "'  -5.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'20.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -4.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Formula = "'155.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Formula = "'  -1.02%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Formula = "'  -3.62%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Formula = "'  -3.61%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Formula = "'  -6.48%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'24.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  -6.56%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Formula = "'  -5.63%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D41").Formula = "'3.019.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Formula = "'3.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  -3.69%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'  -2.81%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Formula = "'2.193.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Formula = "'  -5.78%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'  -5.97%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Formula = "'5.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Formula = "'  -1.76%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Formula = "'0.936"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Formula = "'  -9.51%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Formula = "'0.0235"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Formula = "'  -4.92%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'  -4.34%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Formula = "'  -11.14%  "
$ws.Range("E51").Style = "Normal"
